$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.364.55'
$ws.Range('E2').Value = '  +2.09%  '
$ws.Range('D3').Value = '2.032.85'
$ws.Range('E3').Value = '  +3.21%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '247.20'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.41%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.621'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '58.47'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.28%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.389'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +3.19%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0805'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.34%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.99'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.47%  '
$ws.Range('D13').Value = '2.335.80'
$ws.Range('E13').Value = '  +3.20%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.841'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.74'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.41'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.37%  '
$ws.Range('D17').Value = '2.037.10'
$ws.Range('E17').Value = '  +3.14%  '
$ws.Range('D18').Value = '37.277.76'
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '70.21'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.59%  '
$ws.Range('D20').Value = '0.0₃0857'
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.24'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.31%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '228.83'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.54'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +4.34%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.35'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.26'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.41%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '163.58'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.82%  '
$ws.Range('E28').Value = '  -5.04%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.88'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.77%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.36'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.25%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.121'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0672'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +9.49%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.78'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.52%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.54'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.30%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.50'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +9.88%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.57'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +6.33%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('E38').Value = '  +2.31%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.44'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.98%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.01'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.48%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0973'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.81%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0218'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.15%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.17'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.88%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.57'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.45%  '
$ws.Range('D45').Value = '1.400.40'
$ws.Range('E45').Value = '  +2.55%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '91.39'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.89%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.05'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.49%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.46'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.62%  '
$ws.Range('E49').Value = '  +15.04%  '
$ws.Range('E50').Value = '  +2.15%  '
$ws.Range('D51').Value = '2.226.87'
$ws.Range('E51').Value = '  +3.23%  '
